# Update dSF (column F) values on Sheet1 to reflect repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = -8
    4  = -2
    6  = -5
    7  = -5
    8  = 7
    9  = -6
    10 = -5
    14 = -3
    17 = 3
    18 = -9
    20 = -6
    22 = -4
    23 = -5
    26 = -2
    30 = -2
    31 = 2
    32 = -2
    34 = -2
    35 = -1
    38 = 0
    41 = 2
    43 = 0
    44 = -1
    46 = 4
    49 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}

$wb.Save()
